$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2  = "2025-10-17T07:09:27.081073"
    3  = "2025-10-17T07:09:27.081073"
    4  = "2025-10-17T07:09:27.081073"
    5  = "2025-10-17T07:09:27.081073"
    6  = "2025-10-17T07:09:27.081073"
    7  = "2025-10-17T07:09:27.081073"
    8  = "2025-10-17T07:09:27.081073"
    9  = "2025-10-17T07:09:27.081073"
    10 = "2025-10-17T07:09:27.081073"
    11 = "2025-10-17T07:09:27.081073"
    12 = "2025-10-17T07:09:27.081073"
    13 = "2025-10-17T07:09:27.081073"
    14 = "2025-10-17T07:09:27.081073"
    15 = "2025-10-17T07:09:27.081073"
    16 = "2025-10-17T07:09:27.156081"
    17 = "2025-10-17T07:09:27.156081"
    18 = "2025-10-17T07:09:27.156081"
    19 = "2025-10-17T07:09:27.156081"
    20 = "2025-10-17T07:09:27.156081"
    21 = "2025-10-17T07:09:27.157081"
    22 = "2025-10-17T07:09:27.157081"
    23 = "2025-10-17T07:09:27.157081"
    24 = "2025-10-17T07:09:27.157081"
    25 = "2025-10-17T07:09:27.158084"
    26 = "2025-10-17T07:09:27.210833"
    27 = "2025-10-17T07:09:27.210833"
    28 = "2025-10-17T07:09:27.210833"
    29 = "2025-10-17T07:09:27.210833"
    30 = "2025-10-17T07:09:27.210833"
    31 = "2025-10-17T07:09:27.210833"
    32 = "2025-10-17T07:09:27.210833"
    33 = "2025-10-17T07:09:27.210833"
    34 = "2025-10-17T07:09:27.210833"
    35 = "2025-10-17T07:09:27.210833"
    36 = "2025-10-17T07:09:27.210833"
    37 = "2025-10-17T07:09:27.210833"
    38 = "2025-10-17T07:09:27.210833"
    39 = "2025-10-17T07:09:27.210833"
    40 = "2025-10-17T07:09:27.210833"
    41 = "2025-10-17T07:09:27.210833"
    42 = "2025-10-17T07:09:27.210833"
    43 = "2025-10-17T07:09:27.210833"
    44 = "2025-10-17T07:09:27.210833"
    45 = "2025-10-17T07:09:27.210833"
    46 = "2025-10-17T07:09:27.210833"
    47 = "2025-10-17T07:09:27.210833"
    48 = "2025-10-17T07:09:27.210833"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 26).Value = $timestamps[$row]
}
